$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FromSheet")
$ws2 = $wb.Worksheets.Item("MinValues")
$ws3 = $wb.Worksheets.Item("MaxValues")

# ---------------------------------------------------------------------------
# FromSheet (sheet1): add the "Post 1922" (C) and "Total" (D) columns so the
# split of pre/post-1922 Upper Basin water rights is shown, matching the
# layout already used on the MinValues / MaxValues sheets.
# ---------------------------------------------------------------------------

# Bring in the formatting used for the equivalent header/data cells on the
# MinValues sheet (C3:D7) and MaxValues sheet (D3), so the new cells pick up
# the existing cell styles instead of inventing new ones.
$ws2.Range("C3:D7").Copy() | Out-Null
$ws1.Range("C3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws3.Range("D3").Copy() | Out-Null
$ws1.Range("D3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws1.Range("B4").Copy() | Out-Null
$ws1.Range("C8:D10").PasteSpecial(-4122) | Out-Null

$ws1.Range("B11").Copy() | Out-Null
$ws1.Range("C11:D12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Header row
$ws1.Range("C3").Value = "Post 1922"
$ws1.Range("D3").Value = "Total"

# Utah
$ws1.Range("D4").Value = 0.98
$ws1.Range("C4").Formula = "=D4-B4"

# Wyoming
$ws1.Range("C5").Value = 0.05
$ws1.Range("D5").Formula = "=SUM(B5,C5)"

# Colorado
$ws1.Range("D6").Value = 2.5
$ws1.Range("C6").Formula = "=D6-B6"

# New Mexico
$ws1.Range("D7").Value = 0.415
$ws1.Range("C7").Formula = "=D7-B7"

# Nevada
$ws1.Range("C8").Formula = "=0.3-B8"
$ws1.Range("D8").Formula = "=SUM(B8:C8)"

# Arizona
$ws1.Range("C9").Formula = "=2.8-B9"
$ws1.Range("D9").Formula = "=SUM(B9:C9)"

# California
$ws1.Range("C10").Formula = "=4.4-B10"
$ws1.Range("D10").Formula = "=SUM(B10:C10)"

# Upper Basin total row
$ws1.Range("C11").Formula = "=SUM(C4:C7)"
$ws1.Range("D11").Formula = "=SUM(D4:D7)"

# Lower Basin total row
$ws1.Range("C12").Formula = "=SUM(C8:C10)"
$ws1.Range("D12").Formula = "=SUM(D8:D10)"

# Nevada's "Pre 1922" value (B8) keeps its border style, but now sits amid a
# full row of data -- square up its border to all four sides like the rest
# of the table instead of only the left/right edge it had before.
$ws1.Range("B8").Borders.LineStyle = 1

# Column widths: column C is now a narrow numeric "Post 1922" column instead
# of the wide "Source" column, and column A widened slightly.
$ws1.Columns.Item(1).ColumnWidth = 10.6
$ws1.Columns.Item(3).ColumnWidth = 9.6

# View: scrolled down a bit with the cursor resting past the new columns.
$ws1.Application.ActiveWindow.ScrollRow = 4
$ws1.Range("F10").Select()

# ---------------------------------------------------------------------------
# MinValues: just a change of the selected cell (no data changed).
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C3:D8").Select()

# Re-activate FromSheet so it matches the workbook's tabSelected sheet.
$ws1.Activate()

$wb.Save()
